$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataCell($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.Font.Name = "Calibri"
    $r.Font.Size = 10
    $r.Font.ThemeColor = 1
}

function Set-EmptyStyledCell($addr) {
    $r = $ws.Range($addr)
    $r.Font.Name = "Calibri"
    $r.Font.Size = 10
    $r.Font.ThemeColor = 1
}

# --- Row 2: MCH166-1 ---
Set-DataCell "A2" "MCH166-1"
Set-DataCell "C2" "MISCELLANEOUS, VARIOUS DUTCH BOOKS ON SA"
Set-EmptyStyledCell "D2"
Set-DataCell "E2" "Series"
Set-DataCell "F2" "1 Box"
Set-DataCell "G2" "LOCATION: 22C | GRAP COUNT NUMER: NONE"
Set-EmptyStyledCell "H2"

# --- Row 3: MCH166-2 ---
Set-DataCell "A3" "MCH166-2"
Set-DataCell "C3" "MISCELLANOUS, DUPLICATES"
Set-EmptyStyledCell "D3"
Set-DataCell "E3" "Series"
Set-DataCell "F3" "1 Box"
Set-DataCell "G3" "LOCATION: 22C | GRAP COUNT NUMER: NONE"
Set-EmptyStyledCell "H3"

# --- Row 4: MCH166-3 ---
Set-DataCell "A4" "MCH166-3"
Set-DataCell "C4" "MISCELLANEOUS, VARIOUS DUTCH BOOKS ON SA"
Set-EmptyStyledCell "D4"
Set-DataCell "E4" "Series"
Set-DataCell "F4" "1 Box"
Set-DataCell "G4" "LOCATION: 22C | GRAP COUNT NUMER: NONE"
Set-EmptyStyledCell "H4"

# Restore the frozen header pane + selection as left by the editor.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F16").Select()
